# Apply "Add data for 2022-12-18" changes to the carjacking by month workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet to reflect the new "through" date.
$ws.Name = "Through 2022-12-10"

# Update the December row label text.
$ws.Range("A13").Value = "December (through 12-10)"

# Update December (row 13) values.
$ws.Range("B13").Value = 9
$ws.Range("C13").Value = 28
$ws.Range("D13").Value = 36
$ws.Range("E13").Value = 24
$ws.Range("F13").Value = 14
$ws.Range("G13").Value = 48
$ws.Range("H13").Value = 75
$ws.Range("I13").Value = 44

# Update Total row (row 14) values.
$ws.Range("B14").Value = 300
$ws.Range("C14").Value = 591
$ws.Range("D14").Value = 857
$ws.Range("E14").Value = 706
$ws.Range("F14").Value = 548
$ws.Range("G14").Value = 1312
$ws.Range("H14").Value = 1718
$ws.Range("I14").Value = 1560
